# Apply data + view changes described in the commit.
$wb = $excel.ActiveWorkbook

$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")
$wsInstalacion = $wb.Worksheets.Item("Instalación")

# --- Data changes on "Presupuesto" sheet -------------------------------
# Transporte cost (E24): 15 -> 25
$wsPresupuesto.Range("E24").Value = 25

# Beneficio modulos con bateria (D25): 0.2 -> 0.3
$wsPresupuesto.Range("D25").Value = 0.3

# Beneficio modulos con fuente de poder (D26): 0.2 -> 0.3
$wsPresupuesto.Range("D26").Value = 0.3

# --- View / selection changes -------------------------------------------
# "Instalación" sheet scroll/selection updates (no longer the active tab).
$wsInstalacion.Activate()
$wsInstalacion.Range("I11").Select()
try { $excel.ActiveWindow.ScrollRow = 8 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $wsInstalacion.Range("A8") } catch {}

# "Presupuesto" becomes the active/selected tab, scrolled to row 13,
# with D27 selected.
$wsPresupuesto.Activate()
$wsPresupuesto.Range("D27").Select()
try { $excel.ActiveWindow.ScrollRow = 13 } catch {}
try { $excel.ActiveWindow.TopLeftCell = $wsPresupuesto.Range("A13") } catch {}
